$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.143.54"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.078.54"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.98"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.673"
$ws.Range("E6").Value = "  +2.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.68"
$ws.Range("E7").Value = "  +13.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("E9").Value = "  +4.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "61.61"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  +6.36%  "

$ws.Range("E12").Value = "  +2.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.96"
$ws.Range("E13").Value = "  +5.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.383.94"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.816"
$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("E16").Value = "  +7.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.076.60"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.092.81"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.57"
$ws.Range("E19").Value = "  +10.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "74.57"
$ws.Range("E20").Value = "  +3.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0922"
$ws.Range("E21").Value = "  +9.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.43"
$ws.Range("E22").Value = "  +4.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.02"
$ws.Range("E23").Value = "  -0.52%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -2.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  +13.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.79"
$ws.Range("E27").Value = "  -0.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.29"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("E30").Value = "  +2.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.82"
$ws.Range("E31").Value = "  +7.25%  "

$ws.Range("E32").Value = "  +5.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0634"
$ws.Range("E33").Value = "  +3.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.44"
$ws.Range("E34").Value = "  +8.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0916"
$ws.Range("E35").Value = "  +0.48%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.119"
$ws.Range("E38").Value = "  +29.88%  "

$ws.Range("E39").Value = "  -4.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").Value = "  +1.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.65"
$ws.Range("E41").Value = "  +24.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.97"
$ws.Range("E42").Value = "  -1.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0226"
$ws.Range("E43").Value = "  +0.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.16"
$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.68"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("E46").Value = "  +2.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("E47").Value = "  +12.20%  "

$ws.Range("E48").Value = "  +8.39%  "

$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.303.10"
$ws.Range("E50").Value = "  -1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.94"
$ws.Range("E51").Value = "  -0.36%  "
